$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.871.79"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "2.682.89"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.02"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.79"
$ws.Range("E6").Value = "  -0.66%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  +0.35%  "

$ws.Range("E9").Value = "  +0.42%  "

$ws.Range("E10").Value = "  +2.43%  "

$ws.Range("E11").Value = "  -2.52%  "

$ws.Range("E12").Value = "  +0.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.51"
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000202"
$ws.Range("E14").Value = "  +6.62%  "

$ws.Range("D15").Value = "3.162.20"
$ws.Range("E15").Value = "  +0.75%  "

$ws.Range("D16").Value = "65.648.33"
$ws.Range("E16").Value = "  +0.59%  "

$ws.Range("D17").Value = "2.681.92"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("E18").Value = "  -1.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.84"
$ws.Range("E19").Value = "  -1.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.60"
$ws.Range("E20").Value = "  +3.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.48"
$ws.Range("E21").Value = "  -2.19%  "

$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.49"
$ws.Range("E23").Value = "  +2.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000111"
$ws.Range("E24").Value = "  +7.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.85"
$ws.Range("E25").Value = "  +3.26%  "

$ws.Range("E26").Value = "  -4.71%  "

$ws.Range("E27").Value = "  -1.36%  "

$ws.Range("E28").Value = "  +1.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.18"
$ws.Range("E29").Value = "  -0.96%  "

$ws.Range("E30").Value = "  +0.27%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.17"
$ws.Range("E31").Value = "  -2.14%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "536.49"
$ws.Range("E32").Value = "  -0.87%  "

$ws.Range("E33").Value = "  -3.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.58"
$ws.Range("E34").Value = "  +2.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.42"
$ws.Range("E35").Value = "  -4.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.427"
$ws.Range("E36").Value = "  -1.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.51"
$ws.Range("E37").Value = "  -0.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.54"
$ws.Range("E38").Value = "  -1.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("E40").Value = "  -2.20%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.26"
$ws.Range("E42").Value = "  -0.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "166.73"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.11"
$ws.Range("E44").Value = "  -1.84%  "

$ws.Range("E45").Value = "  +0.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.20"
$ws.Range("E46").Value = "  +0.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.24"
$ws.Range("E47").Value = "  -3.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0264"
$ws.Range("E48").Value = "  -0.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.651"
$ws.Range("E49").Value = "  -1.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.35"
$ws.Range("E50").Value = "  +2.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0990"
$ws.Range("E51").Value = "  +0.36%  "
